# Insert a new row at position 355 (shifts existing rows 355:477 down to 356:478,
# carrying their formatting/values along) and populate the new row with the
# new weekly data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(355).Insert()

$ws.Range("A355").Value = 6
$ws.Range("B355").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C355").Value = "Metropolitana"
$ws.Range("D355").Value = 44588
$ws.Range("E355").Value = 13
$ws.Range("F355").Value = 100112012
$ws.Range("G355").Value = "Espinaca"
$ws.Range("H355").Value = "Sin especificar"
$ws.Range("I355").Value = "Primera"
$ws.Range("J355").Value = 370
$ws.Range("K355").Value = 7500
$ws.Range("L355").Value = 8000
$ws.Range("M355").Value = 7703
$ws.Range("N355").Value = "`$/cuna 10 kilos"
$ws.Range("O355").Value = "Región Metropolitana"
$ws.Range("P355").Value = 770
$ws.Range("Q355").Value = 10
$ws.Range("R355").Value = "Hortaliza"
